# Update 21-Jun-2021, midday update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 32: additional Rp 260,000 wages/meal expense on top of the 60,000 ---
$ws.Range("D32").Formula = "=60000+260000"

# --- Row 33: additional purchases of Rp 25,000,000 + Rp 1,379,000 ---
$ws.Range("D33").Formula = "=775000+577200+1555000+25000000+1379000"

# --- Row 35: additional A/R collections ---
$ws.Range("C35").Formula = "=577200+1000000+30000000+9035000+9375000"

# --- Row 37 (new entry): SALES - cash/retail ---
$ws.Range("B37").Value = "SALES - cash/retail"
$ws.Range("C37").Formula = "=21768725-3041725-9375000"

# --- Row 38 (new entry): SETOR KE BANK ---
$ws.Range("B38").Value = "SETOR KE BANK"
$ws.Range("D38").Value = 30000000

# --- Row 39 (new entry): new day (19-Jun-2021), Wages Expense ---
$ws.Range("A39").Value = 44366
$ws.Range("A39").NumberFormat = "[$-409]d/mmm/yyyy;@"
$ws.Range("B39").Value = "Wages Expense"
$ws.Range("D39").Formula = "=60000"

# --- Row 40 (new entry): BELI plastik ---
$ws.Range("B40").Value = "BELI plastik"
$ws.Range("D40").Formula = "=75000"

# --- Update selection state on the sheet ---
$ws.Range("C37").Select()
